$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("O3").Value = "'5.0"
$ws.Range("Q3").Value = "'7.0"
$ws.Range("S3").Value = "'5.0"
$ws.Range("U3").Value = "'21.0"
$ws.Range("W3").Value = "'5.0"
$ws.Range("O4").Value = "'5.0"
$ws.Range("Q4").Value = "'7.0"
$ws.Range("S4").Value = "'5.0"
$ws.Range("U4").Value = "'21.0"
$ws.Range("W4").Value = "'5.0"
$ws.Range("E5").Value = "'15.0"
$ws.Range("G5").Value = "'3.0"
$ws.Range("H5").Value = "'10.0"
$ws.Range("I5").Value = "'15.0"
$ws.Range("J5").Value = "'3.0"
$ws.Range("K5").Value = "'19.0"
$ws.Range("O5").Value = "'3.0"
$ws.Range("Q5").Value = "'6.0"
$ws.Range("W5").Value = "'3.0"
$ws.Range("C6").Value = "'9.0"
$ws.Range("F6").Value = "'12.0"
$ws.Range("I6").Value = "'11.0"
$ws.Range("K6").Value = "'18.0"
$ws.Range("M6").Value = "'17.0"
$ws.Range("O6").Value = "'10.5"
$ws.Range("Q6").Value = "'10.5"
$ws.Range("S6").Value = "'12.0"
$ws.Range("U6").Value = "'14.5"
$ws.Range("W6").Value = "'11.0"
$ws.Range("E7").Value = "'10.5"
$ws.Range("F7").Value = "'9.0"
$ws.Range("H7").Value = "'2.5"
$ws.Range("I7").Value = "'14.5"
$ws.Range("K7").Value = "'14.0"
$ws.Range("L7").Value = "'7.0"
$ws.Range("M7").Value = "'14.5"
$ws.Range("O7").Value = "'9.0"
$ws.Range("S7").Value = "'9.0"
$ws.Range("U7").Value = "'18.5"
$ws.Range("W7").Value = "'9.0"
$ws.Range("C8").Value = "'16.0"
$ws.Range("E8").Value = "'14.0"
$ws.Range("F8").Value = "'16.0"
$ws.Range("G8").Value = "'15.0"
$ws.Range("H8").Value = "'11.5"
$ws.Range("I8").Value = "'14.5"
$ws.Range("J8").Value = "'15.0"
$ws.Range("L8").Value = "'14.0"
$ws.Range("M8").Value = "'19.0"
$ws.Range("O8").Value = "'15.0"
$ws.Range("Q8").Value = "'15.5"
$ws.Range("S8").Value = "'15.0"
$ws.Range("U8").Value = "'20.0"
$ws.Range("C9").Value = "'12.0"
$ws.Range("E9").Value = "'7.5"
$ws.Range("F9").Value = "'13.0"
$ws.Range("H9").Value = "'13.5"
$ws.Range("K9").Value = "'16.0"
$ws.Range("M9").Value = "'18.0"
$ws.Range("O9").Value = "'12.5"
$ws.Range("Q9").Value = "'11.5"
$ws.Range("S9").Value = "'13.0"
$ws.Range("U9").Value = "'17.5"
$ws.Range("W9").Value = "'12.0"
$ws.Range("C11").Value = "'22.0"
$ws.Range("E11").Value = "'14.5"
$ws.Range("H11").Value = "'17.0"
$ws.Range("I11").Value = "'19.5"
$ws.Range("K11").Value = "'18.0"
$ws.Range("M11").Value = "9.5*"
$ws.Range("O11").Value = "'23.0"
$ws.Range("Q11").Value = "'20.5"
$ws.Range("S11").Value = "'21.0"
$ws.Range("U11").Value = "7.5*"
$ws.Range("W11").Value = "'22.0"
$ws.Range("C14").Value = "'4.5"
$ws.Range("F14").Value = "'3.0"
$ws.Range("G14").Value = "'3.0"
$ws.Range("H14").Value = "'10.0"
$ws.Range("I14").Value = "'16.0"
$ws.Range("J14").Value = "'3.0"
$ws.Range("K14").Value = "'19.5"
$ws.Range("L14").Value = "'4.0"
$ws.Range("O14").Value = "'3.0"
$ws.Range("Q14").Value = "'6.5"
$ws.Range("S14").Value = "'4.5"
$ws.Range("W14").Value = "'4.0"
$ws.Range("C16").Value = "'6.5"
$ws.Range("E16").Value = "'13.0"
$ws.Range("F16").Value = "'7.0"
$ws.Range("G16").Value = "'7.0"
$ws.Range("H16").Value = "'8.5"
$ws.Range("I16").Value = "'14.5"
$ws.Range("J16").Value = "'7.0"
$ws.Range("K16").Value = "'17.0"
$ws.Range("L16").Value = "'7.0"
$ws.Range("M16").Value = "'19.5"
$ws.Range("O16").Value = "'6.5"
$ws.Range("Q16").Value = "'6.5"
$ws.Range("S16").Value = "'6.5"
$ws.Range("U16").Value = "'20.5"
$ws.Range("W16").Value = "'7.0"
$ws.Range("E19").Value = "'11.0"
$ws.Range("F19").Value = "'6.0"
$ws.Range("G19").Value = "'4.0"
$ws.Range("H19").Value = "'7.5"
$ws.Range("I19").Value = "'14.5"
$ws.Range("J19").Value = "'4.0"
$ws.Range("K19").Value = "'17.0"
$ws.Range("L19").Value = "'4.0"
$ws.Range("M19").Value = "'17.5"
$ws.Range("Q19").Value = "'6.0"
$ws.Range("S19").Value = "'6.0"
$ws.Range("U19").Value = "'18.5"
$ws.Range("W19").Value = "'6.0"
